$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 150
$ws.Cells.Item(149, 1).Copy($ws.Cells.Item(150, 1)) | Out-Null
$ws.Cells.Item(149, 5).Copy($ws.Cells.Item(150, 5)) | Out-Null
$ws.Cells.Item(150, 1).Value = 148
$ws.Cells.Item(150, 2).Value = 6788944
$ws.Cells.Item(150, 3).Value = 'Croatia HNL'
$ws.Cells.Item(150, 4).Value = 'Croatia HNL'
$ws.Cells.Item(150, 5).Value = 45396.58333333334
$ws.Cells.Item(150, 6).Value = 'Istra 1961'
$ws.Cells.Item(150, 7).Value = 'HNK Rijeka'
$ws.Cells.Item(150, 8).Value = 0
$ws.Cells.Item(150, 9).Value = 2
$ws.Cells.Item(150, 10).Value = 'A'
$ws.Cells.Item(150, 11).Value = 5.5
$ws.Cells.Item(150, 12).Value = 3.6
$ws.Cells.Item(150, 13).Value = 1.615
$ws.Cells.Item(150, 14).Value = 8.5
$ws.Cells.Item(150, 15).Value = 4
$ws.Cells.Item(150, 16).Value = 1.4
$ws.Cells.Item(150, 17).Value = 1.25
$ws.Cells.Item(150, 18).Value = 1.775
$ws.Cells.Item(150, 19).Value = 2.1
$ws.Cells.Item(150, 20).Value = 2
$ws.Cells.Item(150, 21).Value = 1.8
$ws.Cells.Item(150, 22).Value = 2.05
$ws.Cells.Item(150, 23).Value = -1
$ws.Cells.Item(150, 24).Value = -1
$ws.Cells.Item(150, 25).Value = 0.3999999999999999
$ws.Cells.Item(150, 26).Value = -1
$ws.Cells.Item(150, 27).Value = 1.1
$ws.Cells.Item(150, 28).Value = 0
$ws.Cells.Item(150, 29).Value = 0

# Row 151
$ws.Cells.Item(149, 1).Copy($ws.Cells.Item(151, 1)) | Out-Null
$ws.Cells.Item(149, 5).Copy($ws.Cells.Item(151, 5)) | Out-Null
$ws.Cells.Item(151, 1).Value = 149
$ws.Cells.Item(151, 2).Value = 6962505
$ws.Cells.Item(151, 3).Value = 'Croatia HNL'
$ws.Cells.Item(151, 4).Value = 'Croatia HNL'
$ws.Cells.Item(151, 5).Value = 45402.5
$ws.Cells.Item(151, 6).Value = 'NK Osijek'
$ws.Cells.Item(151, 7).Value = 'Istra 1961'
$ws.Cells.Item(151, 11).Value = 1.55
$ws.Cells.Item(151, 12).Value = 3.8
$ws.Cells.Item(151, 13).Value = 6.5
$ws.Cells.Item(151, 14).Value = 1.6
$ws.Cells.Item(151, 15).Value = 3.75
$ws.Cells.Item(151, 16).Value = 6
$ws.Cells.Item(151, 17).Value = -0.75
$ws.Cells.Item(151, 18).Value = 1.775
$ws.Cells.Item(151, 19).Value = 2.1
$ws.Cells.Item(151, 20).Value = 2.25
$ws.Cells.Item(151, 21).Value = 1.9
$ws.Cells.Item(151, 22).Value = 1.95
$ws.Cells.Item(151, 23).Value = 0
$ws.Cells.Item(151, 24).Value = 0
$ws.Cells.Item(151, 25).Value = 0
$ws.Cells.Item(151, 26).Value = 0
$ws.Cells.Item(151, 27).Value = 0

# Row 152
$ws.Cells.Item(149, 1).Copy($ws.Cells.Item(152, 1)) | Out-Null
$ws.Cells.Item(149, 5).Copy($ws.Cells.Item(152, 5)) | Out-Null
$ws.Cells.Item(152, 1).Value = 150
$ws.Cells.Item(152, 2).Value = 6954858
$ws.Cells.Item(152, 3).Value = 'Croatia HNL'
$ws.Cells.Item(152, 4).Value = 'Croatia HNL'
$ws.Cells.Item(152, 5).Value = 45402.59027777778
$ws.Cells.Item(152, 6).Value = 'Slaven Belupo'
$ws.Cells.Item(152, 7).Value = 'Hajduk Split'
$ws.Cells.Item(152, 11).Value = 5.75
$ws.Cells.Item(152, 12).Value = 4
$ws.Cells.Item(152, 13).Value = 1.571
$ws.Cells.Item(152, 14).Value = 5.5
$ws.Cells.Item(152, 15).Value = 4
$ws.Cells.Item(152, 16).Value = 1.6
$ws.Cells.Item(152, 17).Value = 0.75
$ws.Cells.Item(152, 18).Value = 2.05
$ws.Cells.Item(152, 19).Value = 1.8
$ws.Cells.Item(152, 20).Value = 2.5
$ws.Cells.Item(152, 21).Value = 1.925
$ws.Cells.Item(152, 22).Value = 1.925
$ws.Cells.Item(152, 23).Value = 0
$ws.Cells.Item(152, 24).Value = 0
$ws.Cells.Item(152, 25).Value = 0
$ws.Cells.Item(152, 26).Value = 0
$ws.Cells.Item(152, 27).Value = 0

# Row 153
$ws.Cells.Item(149, 1).Copy($ws.Cells.Item(153, 1)) | Out-Null
$ws.Cells.Item(149, 5).Copy($ws.Cells.Item(153, 5)) | Out-Null
$ws.Cells.Item(153, 1).Value = 151
$ws.Cells.Item(153, 2).Value = 6965778
$ws.Cells.Item(153, 3).Value = 'Croatia HNL'
$ws.Cells.Item(153, 4).Value = 'Croatia HNL'
$ws.Cells.Item(153, 5).Value = 45403.5
$ws.Cells.Item(153, 6).Value = 'NK Lokomotiva Zagreb'
$ws.Cells.Item(153, 7).Value = 'Dinamo Zagreb'
$ws.Cells.Item(153, 11).Value = 6
$ws.Cells.Item(153, 12).Value = 4.333
$ws.Cells.Item(153, 13).Value = 1.5
$ws.Cells.Item(153, 14).Value = 7
$ws.Cells.Item(153, 15).Value = 4.333
$ws.Cells.Item(153, 16).Value = 1.45
$ws.Cells.Item(153, 17).Value = 1.25
$ws.Cells.Item(153, 18).Value = 1.8
$ws.Cells.Item(153, 19).Value = 2.05
$ws.Cells.Item(153, 20).Value = 2.5
$ws.Cells.Item(153, 21).Value = 1.825
$ws.Cells.Item(153, 22).Value = 2.025
$ws.Cells.Item(153, 23).Value = 0
$ws.Cells.Item(153, 24).Value = 0
$ws.Cells.Item(153, 25).Value = 0
$ws.Cells.Item(153, 26).Value = 0
$ws.Cells.Item(153, 27).Value = 0

# Row 154
$ws.Cells.Item(149, 1).Copy($ws.Cells.Item(154, 1)) | Out-Null
$ws.Cells.Item(149, 5).Copy($ws.Cells.Item(154, 5)) | Out-Null
$ws.Cells.Item(154, 1).Value = 152
$ws.Cells.Item(154, 2).Value = 6962506
$ws.Cells.Item(154, 3).Value = 'Croatia HNL'
$ws.Cells.Item(154, 4).Value = 'Croatia HNL'
$ws.Cells.Item(154, 5).Value = 45403.59027777778
$ws.Cells.Item(154, 6).Value = 'HNK Rijeka'
$ws.Cells.Item(154, 7).Value = 'HNK Gorica'
$ws.Cells.Item(154, 11).Value = 1.285
$ws.Cells.Item(154, 12).Value = 6
$ws.Cells.Item(154, 13).Value = 8.5
$ws.Cells.Item(154, 14).Value = 1.25
$ws.Cells.Item(154, 15).Value = 6.5
$ws.Cells.Item(154, 16).Value = 9
$ws.Cells.Item(154, 17).Value = -1.75
$ws.Cells.Item(154, 18).Value = 1.925
$ws.Cells.Item(154, 19).Value = 1.925
$ws.Cells.Item(154, 20).Value = 3
$ws.Cells.Item(154, 21).Value = 1.975
$ws.Cells.Item(154, 22).Value = 1.875
$ws.Cells.Item(154, 23).Value = 0
$ws.Cells.Item(154, 24).Value = 0
$ws.Cells.Item(154, 25).Value = 0
$ws.Cells.Item(154, 26).Value = 0
$ws.Cells.Item(154, 27).Value = 0

# Row 155
$ws.Cells.Item(149, 1).Copy($ws.Cells.Item(155, 1)) | Out-Null
$ws.Cells.Item(149, 5).Copy($ws.Cells.Item(155, 5)) | Out-Null
$ws.Cells.Item(155, 1).Value = 153
$ws.Cells.Item(155, 2).Value = 6957866
$ws.Cells.Item(155, 3).Value = 'Croatia HNL'
$ws.Cells.Item(155, 4).Value = 'Croatia HNL'
$ws.Cells.Item(155, 5).Value = 45404.54166666666
$ws.Cells.Item(155, 6).Value = 'NK Rudes'
$ws.Cells.Item(155, 7).Value = 'NK Varazdin'
$ws.Cells.Item(155, 11).Value = 5
$ws.Cells.Item(155, 12).Value = 3.75
$ws.Cells.Item(155, 13).Value = 1.615
$ws.Cells.Item(155, 14).Value = 4.333
$ws.Cells.Item(155, 15).Value = 3.6
$ws.Cells.Item(155, 16).Value = 1.727
$ws.Cells.Item(155, 17).Value = 0.75
$ws.Cells.Item(155, 18).Value = 1.85
$ws.Cells.Item(155, 19).Value = 2
$ws.Cells.Item(155, 20).Value = 2.5
$ws.Cells.Item(155, 21).Value = 1.925
$ws.Cells.Item(155, 22).Value = 1.925
$ws.Cells.Item(155, 23).Value = 0
$ws.Cells.Item(155, 24).Value = 0
$ws.Cells.Item(155, 25).Value = 0
$ws.Cells.Item(155, 26).Value = 0
$ws.Cells.Item(155, 27).Value = 0
